$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 98.22 -> 0M
$t.Cell(1, 1).Range.Text = "0M"

# Row 2: 15.9 -> 0M
$t.Cell(2, 1).Range.Text = "0M"

# Row 3: 895 -> 0M
$t.Cell(3, 1).Range.Text = "0M"

# Row 4: 20225 -> 20539
$t.Cell(4, 1).Range.Text = "20539"

# Row 7: 0.01441 -> 0.02133
$t.Cell(7, 1).Range.Text = "0.02133"

# Row 8: 0.01234 -> 0.01117
$t.Cell(8, 1).Range.Text = "0.01117"

# Row 9: 0.03669 -> 0.05012
$t.Cell(9, 1).Range.Text = "0.05012"

# Row 10: 0.05473 -> 0.05745
$t.Cell(10, 1).Range.Text = "0.05745"

# Row 11: 0.05673 -> 0.06424
$t.Cell(11, 1).Range.Text = "0.06424"

# Row 12: 4.45280 -> 15.89539
$t.Cell(12, 1).Range.Text = "15.89539"

# Row 44: multi-run tab-separated list -> single value 98.22
$t.Cell(44, 1).Range.Text = "98.22"

# Row 45: multi-run tab-separated list -> single value 15.9
$t.Cell(45, 1).Range.Text = "15.9"

# Row 46: multi-run tab-separated list -> single value 895
$t.Cell(46, 1).Range.Text = "895"
